$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "IsraTv2203.m3u"
$ws.Range("B21").Value = "GroupAnonymousBot"
$ws.Range("C21").Value = 1087968824
$ws.Range("D21").Value = "Playlists"
$ws.Range("E21").Value = "2025-03-22 22:58:56"

$ws.Range("A22").Value = "EGlayList230325.m3u"
$ws.Range("B22").Value = "GroupAnonymousBot"
$ws.Range("C22").Value = 1087968824
$ws.Range("D22").Value = "Playlists"
$ws.Range("E22").Value = "2025-03-23 18:42:03"

$ws.Range("A23").Value = "EGlayList230325.m3u"
$ws.Range("B23").Value = "GroupAnonymousBot"
$ws.Range("C23").Value = 1087968824
$ws.Range("D23").Value = "Playlists"
$ws.Range("E23").Value = "2025-03-23 18:45:17"
